$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 1953.3334
$ws.Range("I46").Value = 1350
$ws.Range("J46").Value = 3160
$ws.Range("K46").Value = 4050
$ws.Range("L46").Value = 9480
$ws.Range("M46").Value = -3931
$ws.Range("N46").Value = -9718

$ws.Range("H60").Value = 1953.3334
$ws.Range("I60").Value = 1350
$ws.Range("J60").Value = 3160
$ws.Range("K60").Value = 4050
$ws.Range("L60").Value = 9480
$ws.Range("M60").Value = -3566
$ws.Range("N60").Value = -10448

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2540.8667
$ws.Range("I2").Value = 1154.125
$ws.Range("J2").Value = 4125.7144
$ws.Range("K2").Value = 1154.125
$ws.Range("L2").Value = 4125.7144
$ws.Range("M2").Value = -1041.125
$ws.Range("N2").Value = -4351.7144

$ws.Range("H32").Value = 395628.8
$ws.Range("I32").Value = 3309.5967
$ws.Range("J32").Value = 2606882.5
$ws.Range("K32").Value = 3309.5967
$ws.Range("L32").Value = 2606882.5
$ws.Range("M32").Value = -3022.5967
$ws.Range("N32").Value = -2607456.5

$ws.Range("H45").Value = 3222.9656
$ws.Range("I45").Value = 3077.842
$ws.Range("J45").Value = 3498.7
$ws.Range("K45").Value = 3077.842
$ws.Range("L45").Value = 3498.7
$ws.Range("M45").Value = -2700.842
$ws.Range("N45").Value = -4252.7

$ws.Range("H61").Value = 2273
$ws.Range("I61").Value = 2008.0625
$ws.Range("K61").Value = 2008.0625
$ws.Range("M61").Value = -1796.0625

$ws.Range("H74").Value = 1359.0605
$ws.Range("I74").Value = 989.8570999999999
$ws.Range("J74").Value = 1631.1052
$ws.Range("K74").Value = 989.8570999999999
$ws.Range("L74").Value = 1631.1052
$ws.Range("M74").Value = -115.8570999999999
$ws.Range("N74").Value = -3379.1052

$ws.Range("H77").Value = 1359.0605
$ws.Range("I77").Value = 989.8570999999999
$ws.Range("J77").Value = 1631.1052
$ws.Range("K77").Value = 4949.2855
$ws.Range("L77").Value = 8155.526
$ws.Range("M77").Value = -581.2855
$ws.Range("N77").Value = -16891.526

$ws.Range("H116").Value = 2540.8667
$ws.Range("I116").Value = 1154.125
$ws.Range("J116").Value = 4125.7144
$ws.Range("K116").Value = 1154.125
$ws.Range("L116").Value = 4125.7144
$ws.Range("M116").Value = 1139.875
$ws.Range("N116").Value = -8713.714400000001

$ws.Range("H122").Value = 15700.349
$ws.Range("I122").Value = 16346.683
$ws.Range("J122").Value = 2450.5
$ws.Range("K122").Value = 49040.049
$ws.Range("L122").Value = 7351.5
$ws.Range("M122").Value = -46590.049
$ws.Range("N122").Value = -12251.5

$ws.Range("H132").Value = 2195.1316
$ws.Range("I132").Value = 1455.5834
$ws.Range("J132").Value = 3462.9285
$ws.Range("K132").Value = 4366.7502
$ws.Range("L132").Value = 10388.7855
$ws.Range("M132").Value = -1836.7502
$ws.Range("N132").Value = -15448.7855

$ws.Range("H136").Value = 2273
$ws.Range("I136").Value = 2008.0625
$ws.Range("K136").Value = 6024.1875
$ws.Range("M136").Value = -3474.1875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2540.8667
$ws.Range("I3").Value = 1154.125
$ws.Range("J3").Value = 4125.7144
$ws.Range("K3").Value = 1154.125
$ws.Range("L3").Value = 4125.7144
$ws.Range("M3").Value = -1040.125
$ws.Range("N3").Value = -4353.7144

$ws.Range("H105").Value = 1598.7858
$ws.Range("I105").Value = 1515.174
$ws.Range("K105").Value = 1515.174
$ws.Range("M105").Value = 231.826

$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws.Range("H134").Value = 6733.44
$ws.Range("I134").Value = 1062.9524
$ws.Range("K134").Value = 3188.857199999999
$ws.Range("M134").Value = -653.8571999999995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1122.2858
$ws.Range("I58").Value = 936.6957
$ws.Range("J58").Value = 1976
$ws.Range("K58").Value = 936.6957
$ws.Range("L58").Value = 1976
$ws.Range("M58").Value = -733.6957
$ws.Range("N58").Value = -2382

$ws.Range("H122").Value = 1528.1522
$ws.Range("I122").Value = 1464.1666
$ws.Range("J122").Value = 1648.125
$ws.Range("K122").Value = 4392.4998
$ws.Range("L122").Value = 4944.375
$ws.Range("M122").Value = -1942.4998
$ws.Range("N122").Value = -9844.375

$ws.Range("H134").Value = 2184.3333
$ws.Range("I134").Value = 1785.3334
$ws.Range("J134").Value = 2583.3333
$ws.Range("K134").Value = 5356.0002
$ws.Range("L134").Value = 7749.999899999999
$ws.Range("M134").Value = -2821.0002
$ws.Range("N134").Value = -12819.9999

$ws.Range("H136").Value = 1122.2858
$ws.Range("I136").Value = 936.6957
$ws.Range("J136").Value = 1976
$ws.Range("K136").Value = 2810.0871
$ws.Range("L136").Value = 5928
$ws.Range("M136").Value = -260.0870999999997
$ws.Range("N136").Value = -11028

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2701.9333
$ws.Range("I102").Value = 2571.7273
$ws.Range("J102").Value = 3060
$ws.Range("K102").Value = 2571.7273
$ws.Range("L102").Value = 3060
$ws.Range("M102").Value = -949.7273
$ws.Range("N102").Value = -6304

$ws.Range("H113").Value = 35715300
$ws.Range("I113").Value = 684.2
$ws.Range("K113").Value = 684.2
$ws.Range("M113").Value = 1485.8

$ws.Range("H122").Value = 2638.7
$ws.Range("I122").Value = 3319.8462
$ws.Range("J122").Value = 1373.7142
$ws.Range("K122").Value = 9959.5386
$ws.Range("L122").Value = 4121.142599999999
$ws.Range("M122").Value = -7509.5386
$ws.Range("N122").Value = -9021.142599999999

$ws.Range("H126").Value = 20835816
$ws.Range("I126").Value = 3992.6667
$ws.Range("J126").Value = 33334910
$ws.Range("K126").Value = 11978.0001
$ws.Range("L126").Value = 100004730
$ws.Range("M126").Value = -9508.000100000001
$ws.Range("N126").Value = -100009670

$ws.Range("H132").Value = 6000.4546
$ws.Range("I132").Value = 7869.8945
$ws.Range("J132").Value = 3463.3572
$ws.Range("K132").Value = 23609.6835
$ws.Range("L132").Value = 10390.0716
$ws.Range("M132").Value = -21079.6835
$ws.Range("N132").Value = -15450.0716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2664.0667
$ws.Range("I7").Value = 2437.9
$ws.Range("J7").Value = 3116.4
$ws.Range("K7").Value = 2437.9
$ws.Range("L7").Value = 3116.4
$ws.Range("M7").Value = -2325.9
$ws.Range("N7").Value = -3340.4

$ws.Range("H40").Value = 2253
$ws.Range("I40").Value = 1981.3478
$ws.Range("J40").Value = 2733.6155
$ws.Range("K40").Value = 1981.3478
$ws.Range("L40").Value = 2733.6155
$ws.Range("M40").Value = -1845.3478
$ws.Range("N40").Value = -3005.6155

$ws.Range("H126").Value = 2664.0667
$ws.Range("I126").Value = 2437.9
$ws.Range("J126").Value = 3116.4
$ws.Range("K126").Value = 7313.700000000001
$ws.Range("L126").Value = 9349.200000000001
$ws.Range("M126").Value = -4843.700000000001
$ws.Range("N126").Value = -14289.2

$ws.Range("H132").Value = 7679.375
$ws.Range("I132").Value = 9499.4
$ws.Range("J132").Value = 4646
$ws.Range("K132").Value = 28498.2
$ws.Range("L132").Value = 13938
$ws.Range("M132").Value = -25968.2
$ws.Range("N132").Value = -18998

$ws.Range("H136").Value = 3845.0476
$ws.Range("I136").Value = 2174.7778
$ws.Range("J136").Value = 13866.667
$ws.Range("K136").Value = 6524.3334
$ws.Range("L136").Value = 41600.001
$ws.Range("M136").Value = -3974.3334
$ws.Range("N136").Value = -46700.001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 23809986
$ws.Range("I113").Value = 382.55554
$ws.Range("K113").Value = 1147.66662
$ws.Range("M113").Value = 1022.33338

$ws.Range("H122").Value = 1661.7894
$ws.Range("I122").Value = 1758.2667
$ws.Range("J122").Value = 1300
$ws.Range("K122").Value = 5274.800099999999
$ws.Range("L122").Value = 3900
$ws.Range("M122").Value = -2824.800099999999
$ws.Range("N122").Value = -8800

$ws.Range("H126").Value = 1872
$ws.Range("I126").Value = 1318.75
$ws.Range("J126").Value = 2676.7273
$ws.Range("K126").Value = 3956.25
$ws.Range("L126").Value = 8030.1819
$ws.Range("M126").Value = -1486.25
$ws.Range("N126").Value = -12970.1819

$ws.Range("H132").Value = 46879810
$ws.Range("I132").Value = 68183260
$ws.Range("J132").Value = 12200.4
$ws.Range("K132").Value = 204549780
$ws.Range("L132").Value = 36601.2
$ws.Range("M132").Value = -204547250
$ws.Range("N132").Value = -41661.2
